# =====================================================================
# 2nd march 2017 commit
#
# "day 1" sheet: extend the daily log through row 38 (tasks for 2nd Mar
# 2017 - adding Spring/Hibernate dependencies, wiring applicationContext
# .xml, seeding the products table/DAO tests) and backfill the two
# trailing cells (D32:H32 / D33:H33) that were left blank on row 31-32
# of the original log.
#
# "error report" sheet: log the ContextComponentScan_Error encountered
# while wiring the Spring context XSD.
# =====================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "day 1"
$ws2 = $wb.Worksheets.Item(2)   # "error report"

# ---------------------------------------------------------------------
# "error report": row 6 was an empty placeholder row; fill it in.
# Column order mirrors how the cells were actually authored (Error ID,
# Full Error, Link, then the short "Solution" placeholder) so that new
# shared-string entries come out in the same order as the source file.
# ---------------------------------------------------------------------
$ws2.Cells.Item(6,1).Value = 'ContextComponentScan_Error'

$ws2.Cells.Item(6,2).Value = 'Multiple annotations found at this line: - schema_reference.4: Failed to read schema document ''http://www.springframework.org/schema/context/spring-context-3.0.xsd'', because 1)   could not find the document; 2) the document could not be read; 3) the root element of the document is not <xsd:schema>. - cvc-complex-type.2.4.c: The matching wildcard is strict, but no declaration can be found for element ''context:component-scan''.'
$ws2.Cells.Item(6,2).WrapText = $true

$ws2.Cells.Item(6,4).Value = 'http://stackoverflow.com/questions/28895990/schema-reference-failed-to-read-schema'
$ws2.Cells.Item(6,4).WrapText = $true

$ws2.Cells.Item(6,3).Value = ' '
$ws2.Cells.Item(6,3).WrapText = $true

# ---------------------------------------------------------------------
# "day 1": rows 32 & 33 already existed (S.No / Date / Task only) -
# backfill the remaining columns, which were always "NA" / time-taken /
# "N" answers for those two tasks.
# ---------------------------------------------------------------------
$ws1.Cells.Item(32,4).Value = 'NA'
$ws1.Cells.Item(32,5).Value = 'NA'
$ws1.Cells.Item(32,6).Value = '30 minutes'
$ws1.Cells.Item(32,7).Value = 'N'
$ws1.Cells.Item(32,8).Value = 'NA'

$ws1.Cells.Item(33,4).Value = 'NA'
$ws1.Cells.Item(33,5).Value = 'NA'
$ws1.Cells.Item(33,6).Value = '40 minutes'
$ws1.Cells.Item(33,7).Value = 'N'
$ws1.Cells.Item(33,8).Value = 'NA'

# ---------------------------------------------------------------------
# "day 1": new rows 34-38, 2nd Mar 2017 entries.
# ---------------------------------------------------------------------
# Row 34
$ws1.Cells.Item(34,1).Value = 33
$ws1.Cells.Item(34,2).Value = '2nd Mar,2017'
$ws1.Cells.Item(34,3).Value = 'adding spring dependency for using hibernate with H2'
$ws1.Cells.Item(34,4).Value = 'NA'
$ws1.Cells.Item(34,5).Value = 'https://mvnrepository.com/artifact/org.springframework/spring-orm/4.2.2.RELEASE https://commons.apache.org/proper/commons-dbcp/ https://mvnrepository.com/artifact/org.springframework/spring-orm/5.0.0.M5'
$ws1.Cells.Item(34,6).Value = '240 minutes'
$ws1.Cells.Item(34,7).Value = 'N'
$ws1.Cells.Item(34,8).Value = 'NA'
$ws1.Rows.Item(34).RowHeight = 57.6

# Row 35
$ws1.Cells.Item(35,1).Value = 34
$ws1.Cells.Item(35,2).Value = '2nd Mar,2017'
$ws1.Cells.Item(35,3).Value = 'adding applicationContext.xml and connect it using web.xml'
$ws1.Cells.Item(35,4).Value = 'NA'
$ws1.Cells.Item(35,5).Value = 'NA'
$ws1.Cells.Item(35,6).Value = '20 minutes'
$ws1.Cells.Item(35,7).Value = 'N'
$ws1.Cells.Item(35,8).Value = 'NA'
$ws1.Rows.Item(35).RowHeight = 28.8

# Row 37 is authored before row 36 in the source file (new shared
# strings appear in that order) - write C37 ahead of C36 to reproduce
# the same shared-string table.
# Row 37
$ws1.Cells.Item(37,1).Value = 36
$ws1.Cells.Item(37,2).Value = '2nd Mar,2017'
$ws1.Cells.Item(37,3).Value = 'Put some sample data in the table and it should be retrieved in the Products Page in the datatable'
$ws1.Cells.Item(37,4).Value = 'NA'
$ws1.Cells.Item(37,5).Value = 'NA'
$ws1.Cells.Item(37,6).Value = '10 minutes'
$ws1.Cells.Item(37,7).Value = 'N'
$ws1.Cells.Item(37,8).Value = 'NA'
$ws1.Rows.Item(37).RowHeight = 43.2

# Row 36
$ws1.Cells.Item(36,1).Value = 35
$ws1.Cells.Item(36,2).Value = '2nd Mar,2017'
$ws1.Cells.Item(36,3).Value = 'Create Beans for DataSource, SessionFactory, TransactionManager in applicationContext.xml file. '
$ws1.Cells.Item(36,4).Value = 'NA'
$ws1.Cells.Item(36,5).Value = 'NA'
$ws1.Cells.Item(36,6).Value = '30 minutes'
$ws1.Cells.Item(36,7).Value = 'N'
$ws1.Cells.Item(36,8).Value = 'NA'
$ws1.Rows.Item(36).RowHeight = 57.6

# Row 38 only has S.No / Date / Task / Time Taken filled in.
$ws1.Cells.Item(38,1).Value = 37
$ws1.Cells.Item(38,2).Value = '2nd Mar,2017'
$ws1.Cells.Item(38,3).Value = 'Write the Junit test cases for checking the DAO classes'
$ws1.Cells.Item(38,6).Value = 'Pending'
$ws1.Rows.Item(38).RowHeight = 28.8

# ---------------------------------------------------------------------
# View state: scroll "day 1" down to the new rows and select B38; leave
# "error report" scrolled to row 4 with A7 selected, then hand focus
# back to "day 1" so it stays the active tab/sheet.
# ---------------------------------------------------------------------
$ws2.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$ws1.Activate()
$ws1.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1

